# Updates cryptos list prices/volumes (GitHub Actions scrape refresh) and
# fixes the Aptos / Binance-Peg-BSC-USD row ordering (rows 27-28 swap).
# D-column values that look numeric are prefixed with a leading apostrophe
# so Excel stores them as text (matching the sheet's original inline-string
# "price" formatting, e.g. "62.796.96") instead of auto-converting to Number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.833.60"
$ws.Range("E2").Value = "  -0.64%  "

$ws.Range("D3").Value = "2.574.30"
$ws.Range("E3").Value = "  +0.43%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'580.75"
$ws.Range("E5").Value = "  -0.66%  "

$ws.Range("D6").Value = "'144.01"
$ws.Range("E6").Value = "  -2.47%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "'0.592"
$ws.Range("E8").Value = "  -1.79%  "

$ws.Range("E9").Value = "  -2.41%  "

$ws.Range("E10").Value = "  -1.56%  "

$ws.Range("E11").Value = "  -0.88%  "

$ws.Range("D12").Value = "'0.349"
$ws.Range("E12").Value = "  -2.92%  "

$ws.Range("D13").Value = "'26.97"
$ws.Range("E13").Value = "  -2.06%  "

$ws.Range("D14").Value = "3.035.22"
$ws.Range("E14").Value = "  +0.49%  "

$ws.Range("D15").Value = "62.735.45"
$ws.Range("E15").Value = "  -0.65%  "

$ws.Range("E16").Value = "  -3.10%  "

$ws.Range("D17").Value = "2.563.35"
$ws.Range("E17").Value = "  -0.40%  "

$ws.Range("D18").Value = "'11.05"
$ws.Range("E18").Value = "  -2.92%  "

$ws.Range("D19").Value = "'340.16"
$ws.Range("E19").Value = "  -0.72%  "

$ws.Range("D20").Value = "'4.32"
$ws.Range("E20").Value = "  -2.62%  "

$ws.Range("D21").Value = "'6.61"
$ws.Range("E21").Value = "  -3.53%  "

$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("D23").Value = "'67.42"
$ws.Range("E23").Value = "  +1.22%  "

$ws.Range("D24").Value = "'1.58"
$ws.Range("E24").Value = "  +7.18%  "

$ws.Range("D25").Value = "'1.60"
$ws.Range("E25").Value = "  -2.89%  "

$ws.Range("E26").Value = "  -3.98%  "

$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").Value = "'7.99"
$ws.Range("E27").Value = "  -1.99%  "

$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("D29").Value = "'8.24"
$ws.Range("E29").Value = "  -3.29%  "

$ws.Range("D30").Value = "'1.92"
$ws.Range("E30").Value = "  -3.89%  "

$ws.Range("D31").Value = "0.0₃0797"
$ws.Range("E31").Value = "  -3.74%  "

$ws.Range("D32").Value = "'453.03"
$ws.Range("E32").Value = "  -3.65%  "

$ws.Range("D34").Value = "'176.08"
$ws.Range("E34").Value = "  -0.37%  "

$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("E36").Value = "  -1.83%  "

$ws.Range("D37").Value = "'18.88"
$ws.Range("E37").Value = "  -2.23%  "

$ws.Range("D38").Value = "'4.46"
$ws.Range("E38").Value = "  -1.34%  "

$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("E40").Value = "  -3.29%  "

$ws.Range("D41").Value = "'39.97"
$ws.Range("E41").Value = "  +0.71%  "

$ws.Range("D42").Value = "'157.36"
$ws.Range("E42").Value = "  +3.92%  "

$ws.Range("E43").Value = "  -4.02%  "

$ws.Range("D44").Value = "'0.634"
$ws.Range("E44").Value = "  +2.76%  "

$ws.Range("D45").Value = "'20.97"
$ws.Range("E45").Value = "  -0.62%  "

$ws.Range("E46").Value = "  -3.34%  "

$ws.Range("E47").Value = "  -2.42%  "

$ws.Range("D48").Value = "'0.0235"
$ws.Range("E48").Value = "  -2.63%  "

$ws.Range("D49").Value = "'17.94"
$ws.Range("E49").Value = "  -3.13%  "

$ws.Range("D50").Value = "'11.40"
$ws.Range("E50").Value = "  +0.16%  "

$ws.Range("E51").Value = "  -4.04%  "
